# Insert a new daily price record for "Feria Lagunitas de Puerto Montt - Zanahoria"
# right before the current row 431. Inserting shifts all the existing rows from
# 431..483 down to 432..484 (dimension grows from A1:R483 to A1:R484), matching
# the target diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(431).Insert()

$ws.Range("A431").Value = 4
$ws.Range("B431").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C431").Value = "Los Lagos"
$ws.Range("D431").Value = 44918
$ws.Range("E431").Value = 10
$ws.Range("F431").Value = 100114013
$ws.Range("G431").Value = "Zanahoria"
$ws.Range("H431").Value = "Sin especificar"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 200
$ws.Range("K431").Value = 13000
$ws.Range("L431").Value = 13000
$ws.Range("M431").Value = 13000
$ws.Range("N431").Value = "$/saco 20 kilos"
$ws.Range("O431").Value = "Chillán"
$ws.Range("P431").Value = 650
$ws.Range("Q431").Value = 20
$ws.Range("R431").Value = "Hortaliza"
